# Updates the cryptos price table (columns D = Price, E = Volume(1h))
# with the latest scraped values. A couple of rows also swap rank
# position (Chainlink/WrappedEther and Hedera/VeChain traded places),
# so their Coin name (B) and Link (C) cells are rewritten too.
#
# Note: several "Price" values (e.g. "1.003", "345.90") look like plain
# numbers to Excel's auto-detection and would otherwise be silently
# converted to numeric cells, losing the original text formatting
# (trailing zeros, etc.) used throughout this sheet. To keep them as
# text - matching the rest of the column - we briefly force the cell's
# NumberFormat to Text ("@") before assigning the value, then clear the
# formatting again so the cell ends up with its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.536.82'
$ws.Range('E2').Value = '  +5.66%  '
$ws.Range('D3').Value = '1.818.99'
$ws.Range('E3').Value = '  +5.94%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '345.90'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3826'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +3.60%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3517'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '49.76'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('E10').Value = '  +4.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07801'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +4.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.29'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +10.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.628'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +6.08%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.254'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.95%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.815.86'
$ws.Range('E16').Value = '  +5.94%  '
$ws.Range('E17').Value = '  +4.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06735'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '86.26'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.70'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +8.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.545'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +7.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.27'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.79%  '
$ws.Range('D24').Value = '27.530.57'
$ws.Range('E24').Value = '  +5.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.459'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('E26').Value = '  +7.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.18'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +15.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.502'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +14.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '153.54'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').Value = '2.020.71'
$ws.Range('E30').Value = '  +6.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '136.55'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.379'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +6.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.085'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '14.06'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +8.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08788'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.701'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.663'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +5.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.7079'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +14.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2282'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +7.09%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06545'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +5.23%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02425'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +6.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.019'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.297'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.73'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6629'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +13.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.033'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +5.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.193'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +9.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '132.84'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07369'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.62'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.58%  '
